$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cursos" (courses) text for the row currently listing
# "FUNDAMENTOS DE LA PROGRAMACION,PROGRAMACION ESTRUCTURADA,PROGRAMACION ORIENTADA A OBJETOS"
# Remove the leading "FUNDAMENTOS DE LA PROGRAMACION," course from the list.
$ws.Range("C8").Value = "PROGRAMACION ESTRUCTURADA,PROGRAMACION ORIENTADA A OBJETOS"

# Reflect the active selection on the sheet at C8, matching the edit location.
$ws.Range("C8").Select()
